$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.539.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.96%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.562.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.14%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -1.54%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'210.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.15%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.490"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -1.62%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'22.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.69%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.250"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.61%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D11").Value = "'0.0867"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.25%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.786.52"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "'1.562.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.05%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  +0.67%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.26%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'27.518.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.84%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'62.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.19%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'224.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +4.05%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  +1.70%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.0₃0706"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.20%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  -1.55%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.00%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'9.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.97%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D25").Value = "'149.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.13%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +2.47%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'6.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.29%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'15.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.65%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -1.25%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  +1.31%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.0470"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.81%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +0.35%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.457.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.52%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'3.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.58%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +3.05%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +0.67%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  -1.58%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.40%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.541"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.07%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +0.65%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'5.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.63%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'2.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.85%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.989"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.58%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +7.72%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.970"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.71%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'64.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.42%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'1.701.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.23%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'86.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.24%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'Cronos"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.0522"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.91%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'BabyDogeCoin"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.0₇0980"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.31%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  -0.88%  "
$ws.Range("E51").Style = "Normal"

